$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Insert a new column before column A to make room for "psnr_dual"
$ws.Columns.Item(1).Insert()

# Header
$ws.Range("A1").Value = "psnr_dual"

# Values for A2:A51
$values = @(27.055125230000002,26.74787808,25.721080220000001,26.228191899999999,25.503093740000001,25.672323309999999,26.12342688,26.187918230000001,27.336064329999999,26.03711294,26.32087035,27.222181379999999,26.628023800000001,26.17912595,25.322383540000001,25.545682769999999,25.330006170000001,26.536697520000001,25.649504189999998,26.68389998,25.545474850000002,25.751300409999999,25.55895379,26.071683499999999,26.13300654,26.081431510000002,25.451833969999999,26.10893583,25.610426830000002,25.733550080000001,26.10088545,25.54653163,26.586508479999999,25.70177425,25.785502309999998,25.871158250000001,26.312546919999999,25.167697740000001,25.7074994,26.054348210000001,25.486215690000002,25.92780248,25.764158340000002,26.162938759999999,25.075447100000002,25.912440530000001,25.613235240000002,25.44511881,26.208813580000001,26.104154510000001)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Fix column widths to cover new column D (was A:C bestFit 12, now A:D)
$ws.Range("A1:D1").Columns.AutoFit()

# Update selection to match target (A2:A51)
$ws.Range("A2:A51").Select()
